$wb = $excel.ActiveWorkbook

# Update the "展览" sheet (exhibition listing)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 3708
$wsExhibit.Range("F6").Value = 39
$wsExhibit.Range("F7").Value = 193

# Update the "全部类型" sheet (all types listing), which duplicates the same events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 3708
$wsAll.Range("F10").Value = 39
$wsAll.Range("F12").Value = 193
